$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.51
$ws.Range("H2").Value = 8
$ws.Range("V2").Value = 1.11
$ws.Range("Y2").Value = 22
$ws.Range("AC2").Value = 970

# Row 3
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 2.36
$ws.Range("AB3").Value = 24
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 13.5
$ws.Range("AF3").Value = 36
$ws.Range("AI3").Value = 27
$ws.Range("AO3").Value = 11

# Row 4
$ws.Range("R4").Value = 1.44

# Row 6
$ws.Range("U6").Value = 2.22
$ws.Range("X6").Value = 970

# Row 7
$ws.Range("N7").Value = 1.02

# Row 8
$ws.Range("F8").Value = 2.42
$ws.Range("G8").Value = 2.78
$ws.Range("H8").Value = 2.7
$ws.Range("I8").Value = 3.3
$ws.Range("K8").Value = 4
$ws.Range("N8").Value = 1.03
$ws.Range("T8").Value = 1.47
$ws.Range("U8").Value = 1.92
$ws.Range("V8").Value = 1.43
$ws.Range("W8").Value = 1.56
$ws.Range("X8").Value = 23
$ws.Range("Y8").Value = 18.5
$ws.Range("AA8").Value = 65
$ws.Range("AD8").Value = 18
$ws.Range("AF8").Value = 25
$ws.Range("AM8").Value = 100

# Row 9
$ws.Range("J9").Value = 2.9
$ws.Range("Y9").Value = 9.800000000000001
$ws.Range("AB9").Value = 6.6

# Row 11
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 3.5
$ws.Range("H11").Value = 2.68
$ws.Range("I11").Value = 3.2
$ws.Range("J11").Value = 2.84
$ws.Range("K11").Value = 3.45
$ws.Range("M11").Value = 1.13
$ws.Range("O11").Value = 1.55
$ws.Range("P11").Value = 1.49
$ws.Range("Q11").Value = 2.56
$ws.Range("R11").Value = 1.17
$ws.Range("S11").Value = 5
$ws.Range("T11").Value = 2.06
$ws.Range("U11").Value = 1.74
$ws.Range("V11").Value = 1.48
$ws.Range("W11").Value = 1.41
$ws.Range("X11").Value = 10.5
$ws.Range("AA11").Value = 60
$ws.Range("AC11").Value = 8
$ws.Range("AE11").Value = 50
$ws.Range("AH11").Value = 970
$ws.Range("AI11").Value = 85
$ws.Range("AJ11").Value = 75
$ws.Range("AK11").Value = 60
$ws.Range("AL11").Value = 90
$ws.Range("AN11").Value = 80
$ws.Range("AO11").Value = 60

# Row 12
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 5.7
$ws.Range("W12").Value = 1.83

# Row 13
$ws.Range("M13").Value = 1.19
$ws.Range("N13").Value = 2.08
$ws.Range("O13").Value = 1.81
$ws.Range("R13").Value = 1.11
$ws.Range("T13").Value = 2.5
$ws.Range("U13").Value = 1.57
$ws.Range("X13").Value = 970
$ws.Range("AA13").Value = 75
$ws.Range("AC13").Value = 970
$ws.Range("AE13").Value = 70
$ws.Range("AF13").Value = 970
$ws.Range("AI13").Value = 140
$ws.Range("AJ13").Value = 65
$ws.Range("AL13").Value = 140
$ws.Range("AN13").Value = 110
$ws.Range("AO13").Value = 130

# Row 14
$ws.Range("O14").Value = 1.29

# Row 15
$ws.Range("F15").Value = 1.97
$ws.Range("I15").Value = 5.1
$ws.Range("J15").Value = 3.2
$ws.Range("K15").Value = 3.55
$ws.Range("N15").Value = 2.92
$ws.Range("O15").Value = 1.45
$ws.Range("P15").Value = 1.63
$ws.Range("Q15").Value = 2.32
$ws.Range("R15").Value = 1.23
$ws.Range("T15").Value = 2
$ws.Range("U15").Value = 1.84
$ws.Range("X15").Value = 12.5
$ws.Range("Y15").Value = 16.5
$ws.Range("Z15").Value = 40
$ws.Range("AA15").Value = 140
$ws.Range("AB15").Value = 8.800000000000001
$ws.Range("AC15").Value = 9.199999999999999
$ws.Range("AD15").Value = 24
$ws.Range("AE15").Value = 90
$ws.Range("AF15").Value = 14
$ws.Range("AG15").Value = 13
$ws.Range("AH15").Value = 27
$ws.Range("AI15").Value = 110
$ws.Range("AJ15").Value = 32
$ws.Range("AK15").Value = 32
$ws.Range("AL15").Value = 60
$ws.Range("AM15").Value = 190
$ws.Range("AN15").Value = 26
$ws.Range("AO15").Value = 130

# Row 16
$ws.Range("F16").Value = 3.6
$ws.Range("G16").Value = 4.1
$ws.Range("H16").Value = 2.4
$ws.Range("I16").Value = 2.6
$ws.Range("J16").Value = 2.84
$ws.Range("K16").Value = 3.1
$ws.Range("M16").Value = 1.15
$ws.Range("N16").Value = 2.36
$ws.Range("O16").Value = 1.62
$ws.Range("R16").Value = 1.15
$ws.Range("S16").Value = 6.2
$ws.Range("T16").Value = 2.24
$ws.Range("U16").Value = 1.69
$ws.Range("V16").Value = 1.62
$ws.Range("W16").Value = 1.33
$ws.Range("X16").Value = 970
$ws.Range("Y16").Value = 970
$ws.Range("Z16").Value = 970
$ws.Range("AA16").Value = 970
$ws.Range("AB16").Value = 970
$ws.Range("AC16").Value = 7.2
$ws.Range("AD16").Value = 970
$ws.Range("AE16").Value = 970
$ws.Range("AJ16").Value = 110
$ws.Range("AK16").Value = 70
$ws.Range("AL16").Value = 130
$ws.Range("AM16").Value = 290
$ws.Range("AN16").Value = 140
$ws.Range("AO16").Value = 970

# Row 17
$ws.Range("I17").Value = 5.5
$ws.Range("J17").Value = 2.18
$ws.Range("N17").Value = 1.03
$ws.Range("O17").Value = 1.01

# Row 18
$ws.Range("G18").Value = 2.18
$ws.Range("J18").Value = 3.1
$ws.Range("M18").Value = 1.1
$ws.Range("N18").Value = 2.84
$ws.Range("O18").Value = 1.46
$ws.Range("P18").Value = 1.6
$ws.Range("R18").Value = 1.22
$ws.Range("S18").Value = 4.4
$ws.Range("T18").Value = 2.02
$ws.Range("U18").Value = 1.8
$ws.Range("W18").Value = 1.84
$ws.Range("X18").Value = 10.5
$ws.Range("Y18").Value = 13.5
$ws.Range("Z18").Value = 34
$ws.Range("AA18").Value = 140
$ws.Range("AB18").Value = 7.6
$ws.Range("AC18").Value = 7.6
$ws.Range("AD18").Value = 21
$ws.Range("AE18").Value = 75
$ws.Range("AF18").Value = 12.5
$ws.Range("AG18").Value = 12
$ws.Range("AH18").Value = 23
$ws.Range("AI18").Value = 110
$ws.Range("AJ18").Value = 29
$ws.Range("AK18").Value = 29
$ws.Range("AL18").Value = 55
$ws.Range("AM18").Value = 190
$ws.Range("AN18").Value = 29
$ws.Range("AO18").Value = 120
